# Auto-generated Excel COM-interop script to apply cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value. D-column (Price) values are forced to Text
# number format BEFORE assignment so that numeric-looking strings (e.g. '0.999',
# '10.16') are preserved exactly as text instead of being auto-converted to numbers,
# matching the original inlineStr/text cell representation.
$updates = [ordered]@{
    'D2' = '51.188.59'
    'E2' = '  -1.37%  '
    'D3' = '2.769.26'
    'E3' = '  -0.26%  '
    'E4' = '  +0.04%  '
    'D5' = '353.34'
    'E5' = '  -0.32%  '
    'D6' = '107.41'
    'E6' = '  -1.24%  '
    'D7' = '0.548'
    'E7' = '  -2.30%  '
    'D8' = '0.999'
    'E8' = '  +0.01%  '
    'E9' = '  -1.01%  '
    'D10' = '39.40'
    'E10' = '  -1.61%  '
    'E11' = '  +3.40%  '
    'D12' = '0.0830'
    'E12' = '  -2.23%  '
    'D13' = '19.95'
    'E13' = '  +3.17%  '
    'E14' = '  -0.77%  '
    'D15' = '3.202.89'
    'E15' = '  -0.15%  '
    'D16' = '2.764.01'
    'E16' = '  -0.03%  '
    'D17' = '0.926'
    'E17' = '  -0.17%  '
    'D18' = '51.118.70'
    'E18' = '  -1.24%  '
    'E19' = '  +3.61%  '
    'E20' = '  -1.55%  '
    'D21' = '13.07'
    'E21' = '  +0.48%  '
    'D22' = '0.0₃0959'
    'E22' = '  -1.33%  '
    'D23' = '69.58'
    'E23' = '  -0.06%  '
    'D24' = '265.62'
    'E24' = '  -3.07%  '
    'D25' = '2.71'
    'E25' = '  -0.28%  '
    'E26' = '  +0.05%  '
    'D27' = '25.93'
    'E27' = '  -2.16%  '
    'D28' = '0.162'
    'E28' = '  +12.81%  '
    'B29' = 'Toncoin'
    'C29' = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
    'D29' = '2.30'
    'E29' = '  +3.33%  '
    'B30' = 'Cosmos'
    'C30' = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    'D30' = '10.16'
    'E30' = '  +0.49%  '
    'D31' = '35.02'
    'E31' = '  +4.22%  '
    'D32' = '51.80'
    'E32' = '  +0.86%  '
    'D33' = '6.05'
    'E33' = '  +6.38%  '
    'B34' = 'VeChain'
    'C34' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'D34' = '0.0442'
    'E34' = '  -4.48%  '
    'B35' = 'RenderToken'
    'C35' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D35' = '5.51'
    'E35' = '  +3.80%  '
    'D36' = '0.0824'
    'E36' = '  -1.14%  '
    'E37' = '  +0.02%  '
    'D38' = '18.14'
    'E38' = '  +0.53%  '
    'E39' = '  -1.82%  '
    'D40' = '1.96'
    'E40' = '  -1.62%  '
    'B41' = 'Stellar'
    'C41' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'D41' = '0.114'
    'E41' = '  -0.30%  '
    'B42' = 'Stacks'
    'C42' = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
    'D42' = '2.50'
    'E42' = '  -1.02%  '
    'D43' = '120.63'
    'E43' = '  -0.50%  '
    'D44' = '21.95'
    'E44' = '  +0.50%  '
    'E45' = '  -2.50%  '
    'D46' = '2.092.74'
    'E46' = '  +2.00%  '
    'E47' = '  -0.28%  '
    'D48' = '2.27'
    'E48' = '  +0.27%  '
    'D49' = '0.910'
    'E49' = '  -1.26%  '
    'E51' = '  +7.39%  '
}

foreach ($cellRef in $updates.Keys) {
    $col = ($cellRef -replace '[0-9]+$', '')
    $range = $ws.Range($cellRef)
    if ($col -eq 'D') {
        # Ensure text storage so numeric-looking price strings keep their exact
        # formatting (trailing zeros, thousand separators used as decimal dots, etc.)
        $range.NumberFormat = '@'
    }
    $range.Value = $updates[$cellRef]
}
